$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update PLC live data values (commit: Update PLC data 2025-10-13 13:42:11)
$ws.Range("C3").Value = 155010
$ws.Range("C4").Value = 146134
$ws.Range("C5").Value = 8876
$ws.Range("C8").Value = 63.67
